$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'34.392.57"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -0.80%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.802.53"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.35%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.34%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'227.51"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'0.579"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +3.80%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.28%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E9').Value = "'  +0.38%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'  -0.82%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.0953"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +0.29%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'2.062.25"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +0.22%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'11.17"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +0.05%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'1.796.52"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +0.10%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('E15').Value = "'  +0.07%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'34.377.98"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -0.61%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.87%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'68.98"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +0.05%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('B19').Value = "'BitcoinCash"
$ws.Range('B19').Style = 'Normal'
$ws.Range('C19').Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range('C19').Style = 'Normal'
$ws.Range('D19').Value = "'244.76"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -1.47%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('B20').Value = "'ShibaInu"
$ws.Range('B20').Style = 'Normal'
$ws.Range('C20').Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range('C20').Style = 'Normal'
$ws.Range('D20').Value = "'0.0₃0794"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -2.26%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'11.49"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +1.30%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'  +0.38%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = "'  -0.95%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'170.69"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +3.30%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'2.11"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +2.48%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  +3.74%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'16.71"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +0.91%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  +1.66%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  +0.13%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'3.96"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -0.57%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('B31').Value = "'Hedera"
$ws.Range('B31').Style = 'Normal'
$ws.Range('C31').Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range('C31').Style = 'Normal'
$ws.Range('D31').Value = "'0.0527"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +0.80%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('B32').Value = "'PancakeSwap"
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = "'1.24"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +0.81%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  +0.06%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  -0.21%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'1.397.30"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -1.76%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'0.676"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +0.09%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'2.50"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -2.88%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  +0.18%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  -1.97%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'82.87"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -3.17%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'  +2.50%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'0.945"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +0.91%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'2.38"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -0.39%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'13.55"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -0.19%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'  +2.74%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.0512"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -2.25%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'5.98"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -1.67%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'1.962.55"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +0.24%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'104.33"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -1.77%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  +0.33%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = "'  -2.48%  "
$ws.Range('E51').Style = 'Normal'
